# Apply crypto list updates (prices/volumes refreshed, row 35/36 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.444.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.798.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.17'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.601'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '38.94'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +6.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.287'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0667'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.90%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.058.30'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.86'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.791.08'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.394.26'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.626'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.35'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.93'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '238.73'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0762'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.27%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.93%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.08'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.18'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.56'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.64'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.121'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.73'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0512'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.83'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.62%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.72%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.637'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.25%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.300.42'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -7.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0185'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.57%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.20%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.42'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.940'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.00'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0518'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.958.64'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.73%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.66'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.51%  '
